# Updated cryptos list — refresh Price (column D) and Volume(1h) (column E)
# figures for each coin row, and correct the sort order of the
# Chainlink / WrappedBTC rows (17 and 18), per the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price/Volume cell is stored as plain text in the sheet (values such as
# "69.910.16" or "0.998" and "  +0.15%  "). Assigning a plain numeric-looking
# string straight to .Value lets Excel auto-convert it to a real number, which
# would change the cell's type/formatting and break the "69.910.16"-style
# thousands-dot strings. To keep every updated cell as literal text (same as
# before the edit), write it as a quote-prefixed formula (forces text) and
# then reapply the "Normal" cell style so no stray number-format style sticks
# around on the cell.
function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $range.Formula = "'" + $value
    $range.Style = "Normal"
}

$updates = @(
    @{ Cell = 'D2'; Value = '69.910.16' },
    @{ Cell = 'E2'; Value = '  +0.15%  ' },
    @{ Cell = 'D3'; Value = '3.540.63' },
    @{ Cell = 'E3'; Value = '  +1.31%  ' },
    @{ Cell = 'D4'; Value = '0.998' },
    @{ Cell = 'E4'; Value = '  -0.02%  ' },
    @{ Cell = 'D5'; Value = '605.62' },
    @{ Cell = 'E5'; Value = '  -0.29%  ' },
    @{ Cell = 'D6'; Value = '194.72' },
    @{ Cell = 'E6'; Value = '  -0.38%  ' },
    @{ Cell = 'E7'; Value = '  -0.37%  ' },
    @{ Cell = 'E8'; Value = '  +0.05%  ' },
    @{ Cell = 'E9'; Value = '  -3.38%  ' },
    @{ Cell = 'D10'; Value = '0.648' },
    @{ Cell = 'E10'; Value = '  -1.46%  ' },
    @{ Cell = 'D11'; Value = '53.41' },
    @{ Cell = 'E11'; Value = '  -0.43%  ' },
    @{ Cell = 'D12'; Value = '0.0000304' },
    @{ Cell = 'E12'; Value = '  -0.95%  ' },
    @{ Cell = 'E13'; Value = '  -1.00%  ' },
    @{ Cell = 'D14'; Value = '4.096.10' },
    @{ Cell = 'E14'; Value = '  +0.91%  ' },
    @{ Cell = 'D15'; Value = '592.77' },
    @{ Cell = 'E15'; Value = '  -1.44%  ' },
    @{ Cell = 'D16'; Value = '12.79' },
    @{ Cell = 'E16'; Value = '  +1.13%  ' },
    @{ Cell = 'D19'; Value = '3.535.57' },
    @{ Cell = 'E19'; Value = '  +1.18%  ' },
    @{ Cell = 'E20'; Value = '  +1.80%  ' },
    @{ Cell = 'D21'; Value = '0.985' },
    @{ Cell = 'E21'; Value = '  -0.49%  ' },
    @{ Cell = 'D22'; Value = '17.80' },
    @{ Cell = 'E22'; Value = '  -0.89%  ' },
    @{ Cell = 'D23'; Value = '103.16' },
    @{ Cell = 'E23'; Value = '  -1.58%  ' },
    @{ Cell = 'D24'; Value = '5.14' },
    @{ Cell = 'E24'; Value = '  +1.23%  ' },
    @{ Cell = 'D25'; Value = '4.64' },
    @{ Cell = 'E25'; Value = '  -0.14%  ' },
    @{ Cell = 'D26'; Value = '3.05' },
    @{ Cell = 'E26'; Value = '  -1.15%  ' },
    @{ Cell = 'D27'; Value = '10.74' },
    @{ Cell = 'E27'; Value = '  -1.85%  ' },
    @{ Cell = 'D28'; Value = '9.53' },
    @{ Cell = 'E28'; Value = '  -2.85%  ' },
    @{ Cell = 'D29'; Value = '33.19' },
    @{ Cell = 'E29'; Value = '  -2.16%  ' },
    @{ Cell = 'D30'; Value = '7.05' },
    @{ Cell = 'E30'; Value = '  -2.26%  ' },
    @{ Cell = 'D31'; Value = '4.24' },
    @{ Cell = 'E31'; Value = '  -5.35%  ' },
    @{ Cell = 'D32'; Value = '12.33' },
    @{ Cell = 'D33'; Value = '0.116' },
    @{ Cell = 'E33'; Value = '  -0.09%  ' },
    @{ Cell = 'D34'; Value = '63.44' },
    @{ Cell = 'E34'; Value = '  -1.00%  ' },
    @{ Cell = 'D35'; Value = '3.24' },
    @{ Cell = 'E35'; Value = '  +6.72%  ' },
    @{ Cell = 'D36'; Value = '3.834.80' },
    @{ Cell = 'E36'; Value = '  +4.22%  ' },
    @{ Cell = 'D37'; Value = '0.0₃0824' },
    @{ Cell = 'E37'; Value = '  +4.27%  ' },
    @{ Cell = 'E38'; Value = '  +0.28%  ' },
    @{ Cell = 'D39'; Value = '510.91' },
    @{ Cell = 'E39'; Value = '  -1.69%  ' },
    @{ Cell = 'E40'; Value = '  +0.01%  ' },
    @{ Cell = 'D41'; Value = '3.57' },
    @{ Cell = 'E41'; Value = '  -0.11%  ' },
    @{ Cell = 'E42'; Value = '  -0.74%  ' },
    @{ Cell = 'E43'; Value = '  -2.68%  ' },
    @{ Cell = 'D44'; Value = '0.0448' },
    @{ Cell = 'E44'; Value = '  -2.62%  ' },
    @{ Cell = 'D45'; Value = '3.33' },
    @{ Cell = 'E45'; Value = '  +0.46%  ' },
    @{ Cell = 'D46'; Value = '2.82' },
    @{ Cell = 'E46'; Value = '  -1.14%  ' },
    @{ Cell = 'E47'; Value = '  -0.73%  ' },
    @{ Cell = 'E48'; Value = '  +0.00%  ' },
    @{ Cell = 'D49'; Value = '8.51' },
    @{ Cell = 'E49'; Value = '  -2.72%  ' },
    @{ Cell = 'D50'; Value = '0.000246' },
    @{ Cell = 'E50'; Value = '  +3.08%  ' },
    @{ Cell = 'E51'; Value = '  +2.18%  ' }
)

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}

# Rows 17/18: Chainlink and WrappedBTC swapped places in the ranking; update
# the coin name, link, price and volume for both rows accordingly.
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D17' '19.11'
$ws.Range('E17').Value = '  +0.89%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '69.960.50'
$ws.Range('E18').Value = '  +0.08%  '
